# Generate Report for Handoff
# Adds a new handoff record (84860e38-...) as row 3 on each of the three
# sheets: Overview, zh-cn, de-de. Mirrors the existing b9bb5299-... row.

$wb = $excel.ActiveWorkbook

$guidStem  = "84860e38-d757-49f7-827f-af046086e98d"
$oPad126   = "ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo"
$fileName  = $guidStem + $oPad126 + ".md"
$filePath  = "e2e\" + $fileName

$hoStem    = "84860e38-d757-49f7-827f-af046086e98d"
$oPad44    = "oooooooooooooooooooooooooooooooooooooooooo"
$hashPart  = "5c9b098f3c43f4e65ffbc1f4e2b519d4e75001bd"
$xlfZhCn   = $hoStem + $oPad44 + "." + $hashPart + ".zh-cn.xlf"
$xlfDeDe   = $hoStem + $oPad44 + "." + $hashPart + ".de-de.xlf"

$statusText   = "Ready for handoff"
$handoffDate  = "2016-08-30 02:29:52"
$handoffDateZh = "2016-08-30 02:29:48"
$commitHash   = "e5d0ee46a0321d3ba27e9b14ba300e87133f65be"
$githubUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/" + $commitHash + "/e2e/" + $fileName

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview" (table3 / columns A-G)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $fileName
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Range("G3").Value = $handoffDate
$wsOverview.Range("G3").NumberFormat = $dateFormat

$wsOverview.Range("B3").Value = $filePath
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $githubUrl, "", "", $filePath) | Out-Null
$wsOverview.Range("B3").Style = "Hyperlink"
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("B3").Font.Color = 15570276

# ---------------------------------------------------------------------
# Sheet "zh-cn" (table1 / columns A-P)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $fileName
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $xlfZhCn
$wsZhCn.Range("H3").Value = $handoffDateZh
$wsZhCn.Range("H3").NumberFormat = $dateFormat
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = $dateFormat
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $githubUrl, "", "", $fileName) | Out-Null
$wsZhCn.Range("A3").Style = "Hyperlink"
$wsZhCn.Range("A3").Font.Underline = $true
$wsZhCn.Range("A3").Font.Color = 15570276

# ---------------------------------------------------------------------
# Sheet "de-de" (table2 / columns A-P)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $fileName
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $xlfDeDe
$wsDeDe.Range("H3").Value = $handoffDate
$wsDeDe.Range("H3").NumberFormat = $dateFormat
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = $dateFormat
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $githubUrl, "", "", $fileName) | Out-Null
$wsDeDe.Range("A3").Style = "Hyperlink"
$wsDeDe.Range("A3").Font.Underline = $true
$wsDeDe.Range("A3").Font.Color = 15570276

# ---------------------------------------------------------------------
# Column width refresh (E/F on Overview, C on zh-cn/de-de widen slightly
# to accommodate the new "Ready for handoff" text).
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.38
$wsOverview.Columns.Item(6).ColumnWidth = 16.38
$wsZhCn.Columns.Item(3).ColumnWidth = 16.38
$wsDeDe.Columns.Item(3).ColumnWidth = 16.38

Write-Output "Generate Report for Handoff: added 84860e38 handback row to Overview, zh-cn, de-de."
